$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin'
$ws.Range("D2").Value = '91.698.57'
$ws.Range("E2").Value = '  +0.96%  '

# Row 3: 'Ethereum'
$ws.Range("D3").Value = '3.122.49'
$ws.Range("E3").Value = '  +0.17%  '

# Row 4: 'TetherUSD'
$ws.Range("E4").Value = '  +0.08%  '

# Row 5: 'Solana'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.56'
$ws.Range("E5").Value = '  -0.18%  '

# Row 6: 'BNB'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '619.10'
$ws.Range("E6").Value = '  -0.93%  '

# Row 7: 'XRP'
$ws.Range("E7").Value = '  -5.19%  '

# Row 8: 'Dogecoin'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.388'
$ws.Range("E8").Value = '  +4.67%  '

# Row 9: 'USDC'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  -0.01%  '

# Row 10: 'LidoStakedEther'
$ws.Range("D10").Value = '3.122.35'
$ws.Range("E10").Value = '  +0.28%  '

# Row 11: 'Cardano'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.760'
$ws.Range("E11").Value = '  -1.12%  '

# Row 12: 'TRON'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.205'
$ws.Range("E12").Value = '  +0.04%  '

# Row 13: 'ShibaInu'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000254'
$ws.Range("E13").Value = '  +0.79%  '

# Row 14: 'Avalanche'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.37'
$ws.Range("E14").Value = '  -0.13%  '

# Row 15: 'Toncoin'
$ws.Range("E15").Value = '  +1.98%  '

# Row 16: 'WrappedBTC'
$ws.Range("D16").Value = '91.372.80'
$ws.Range("E16").Value = '  +0.89%  '

# Row 17: 'WrappedliquidstakedEther2.0'
$ws.Range("D17").Value = '3.707.51'
$ws.Range("E17").Value = '  +0.62%  '

# Row 18: 'WrappedEther'
$ws.Range("D18").Value = '3.128.16'
$ws.Range("E18").Value = '  +1.23%  '

# Row 19: 'SuiNetwork'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.80'
$ws.Range("E19").Value = '  +0.57%  '

# Row 20: 'Chainlink'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.97'
$ws.Range("E20").Value = '  +3.39%  '

# Row 21: 'Polkadot'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.94'
$ws.Range("E21").Value = '  +0.24%  '

# Row 22: 'BitcoinCash'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '459.06'
$ws.Range("E22").Value = '  +1.66%  '

# Row 23: 'PEPE'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.0000203'
$ws.Range("E23").Value = '  -3.59%  '

# Row 24: 'Uniswap'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.24'
$ws.Range("E24").Value = '  +1.58%  '

# Row 25: 'NEARProtocol'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.94'
$ws.Range("E25").Value = '  +4.76%  '

# Row 26: 'Litecoin'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '89.49'
$ws.Range("E26").Value = '  -4.47%  '

# Row 27: 'Binance-PegBSC-USD'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.54'
$ws.Range("E27").Value = '  +54.62%  '

# Row 28: 'Aptos'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.75'
$ws.Range("E28").Value = '  -1.33%  '

# Row 29: 'Hedera' -> 'WrappedeETH'
$ws.Range("B29").Value = 'WrappedeETH'
$ws.Range("C29").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D29").Value = '3.291.92'
$ws.Range("E29").Value = '  +0.66%  '

# Row 30: 'WrappedeETH' -> 'Hedera'
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.148'
$ws.Range("E30").Value = '  +25.04%  '

# Row 31: 'Dai'
$ws.Range("E31").Value = '  -0.11%  '

# Row 32: 'Stellar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.229'
$ws.Range("E32").Value = '  -1.41%  '

# Row 33: 'Cronos'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.168'
$ws.Range("E33").Value = '  -6.88%  '

# Row 34: 'InternetComputer(DFINITY)'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '9.36'
$ws.Range("E34").Value = '  +2.72%  '

# Row 35: 'Kaspa'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.175'
$ws.Range("E35").Value = '  +7.14%  '

# Row 36: 'EthereumClassic'
$ws.Range("E36").Value = '  -1.53%  '

# Row 37: 'RenderToken'
$ws.Range("E37").Value = '  -3.28%  '

# Row 38: 'PancakeSwap'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.97'
$ws.Range("E38").Value = '  +2.24%  '

# Row 39: 'MantraDAO' -> 'Bittensor'
$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '491.95'
$ws.Range("E39").Value = '  -0.93%  '

# Row 40: 'Bittensor' -> 'MantraDAO'
$ws.Range("B40").Value = 'MantraDAO'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.90'
$ws.Range("E40").Value = '  -7.79%  '

# Row 41: 'Fetch.AI'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.31'
$ws.Range("E41").Value = '  +1.20%  '

# Row 42: 'PolygonEcosystemToken'
$ws.Range("E42").Value = '  +4.27%  '

# Row 43: 'dogwifhat'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.41'
$ws.Range("E43").Value = '  -5.12%  '

# Row 44: 'WhiteBITCoin'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.16'
$ws.Range("E44").Value = '  +0.13%  '

# Row 45: 'USDe'
$ws.Range("E45").Value = '  -0.06%  '

# Row 46: 'ARBITRUM' -> 'Monero'
$ws.Range("B46").Value = 'Monero'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '157.73'
$ws.Range("E46").Value = '  +0.09%  '

# Row 47: 'Stacks' -> 'ARBITRUM'
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.707'

# Row 48: 'Monero' -> 'Stacks'
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.93'
$ws.Range("E48").Value = '  +0.43%  '

# Row 49: 'ImmutableX'
$ws.Range("E49").Value = '  +0.30%  '

# Row 50: 'Filecoin'
$ws.Range("E50").Value = '  -1.66%  '

# Row 51: 'VeChain'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0328'
$ws.Range("E51").Value = '  +3.46%  '
